# Atualização de bases das ligas, do dia: 21-02-2024 às 23:25
#
# The source rows 147-154 (all sharing the same fixture date) have their
# match data (columns B..AC) redistributed among the 8 physical rows.
# Column A (sequential index) and column E (shared date) stay put; we only
# move id/teams/odds (B, F..AC). Row 153 keeps its own data (maps to itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the B:AC payload of every source row BEFORE any writes,
#    so overlapping reads/writes don't clobber data we still need.
$row147 = $ws.Range("B147:AC147").Value2
$row148 = $ws.Range("B148:AC148").Value2
$row149 = $ws.Range("B149:AC149").Value2
$row150 = $ws.Range("B150:AC150").Value2
$row151 = $ws.Range("B151:AC151").Value2
$row152 = $ws.Range("B152:AC152").Value2
$row153 = $ws.Range("B153:AC153").Value2
$row154 = $ws.Range("B154:AC154").Value2

# 2) Write each snapshot back out to its new physical row.
$ws.Range("B147:AC147").Value2 = $row150
$ws.Range("B148:AC148").Value2 = $row149
$ws.Range("B149:AC149").Value2 = $row154
$ws.Range("B150:AC150").Value2 = $row151
$ws.Range("B151:AC151").Value2 = $row147
$ws.Range("B152:AC152").Value2 = $row148
$ws.Range("B153:AC153").Value2 = $row153
$ws.Range("B154:AC154").Value2 = $row152

# 3) A handful of standalone odds tweaks elsewhere in the sheet.
$ws.Range("Q341").Value2 = -0.5
$ws.Range("R341").Value2 = 2.1
$ws.Range("S341").Value2 = 1.775

$ws.Range("R342").Value2 = 1.975
$ws.Range("S342").Value2 = 1.875
$ws.Range("U342").Value2 = 1.9
$ws.Range("V342").Value2 = 1.95

$ws.Range("R345").Value2 = 1.925
$ws.Range("S345").Value2 = 1.925

$ws.Range("R349").Value2 = 1.775
$ws.Range("S349").Value2 = 2.1
